$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.949.50'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.858.01'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.67'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5118'
$ws.Range('E7').Value = '  +2.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3811'
$ws.Range('E8').Value = '  -1.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08270'
$ws.Range('E9').Value = '  -9.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.66'
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.106'
$ws.Range('E11').Value = '  -1.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.178'
$ws.Range('E12').Value = '  -2.37%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.40'
$ws.Range('E13').Value = '  -1.58%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.856.34'
$ws.Range('E14').Value = '  -1.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.180'
$ws.Range('E15').Value = '  -1.36%  '
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001093'
$ws.Range('E17').Value = '  -1.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.19'
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06600'
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.68'
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.998'
$ws.Range('E22').Value = '  -2.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.979.48'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.99'
$ws.Range('E24').Value = '  -3.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.222'
$ws.Range('E25').Value = '  -3.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.574'
$ws.Range('E26').Value = '  +1.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.070.80'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '156.79'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.38'
$ws.Range('E29').Value = '  -1.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.18'
$ws.Range('E30').Value = '  -1.94%  '
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.037'
$ws.Range('E32').Value = '  -2.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.583'
$ws.Range('E33').Value = '  -0.13%  '
$ws.Range('E34').Value = '  +0.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.502'
$ws.Range('E35').Value = '  +1.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06501'
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02402'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2155'
$ws.Range('E38').Value = '  -1.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.203'
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6434'
$ws.Range('E40').Value = '  +0.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.228'
$ws.Range('E41').Value = '  -4.58%  '
$ws.Range('E42').Value = '  -2.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.852'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6118'
$ws.Range('E44').Value = '  +1.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.00'
$ws.Range('E45').Value = '  -2.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.281'
$ws.Range('E46').Value = '  -0.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.657'
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.970'
$ws.Range('E48').Value = '  -1.02%  '
$ws.Range('E49').Value = '  -0.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '120.29'
$ws.Range('E50').Value = '  -0.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.68'
$ws.Range('E51').Value = '  +1.24%  '
